{"js": "// feat: New replacement variables OFFICE, COMPANY\n//\n// Replace the hard-coded company name in the signature block with the\n// new $CURRENTUSERCOMPANY$ replacement variable, and drop the stale\n// \"_GoBack\" bookmark left over from the last edit position (Word\n// normally discards this bookmark the next time the document is saved\n// after being edited).\n\nconst body = context.document.body;\n\n// 1) Swap the literal company name for the new placeholder variable.\nconst results = body.search(\"Super Duper Inc.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"$CURRENTUSERCOMPANY$\", Word.InsertLocation.replace);\n}\n\n// 2) Remove the leftover \"_GoBack\" bookmark, if present.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# feat: New replacement variables OFFICE, COMPANY\n#\n# Replace the hard-coded company name in the signature block with the\n# new $CURRENTUSERCOMPANY$ replacement variable, and drop the stale\n# \"_GoBack\" bookmark left over from the last edit position (Word\n# normally discards this bookmark the next time the document is saved\n# after being edited).\n\n$d = $word.ActiveDocument\n\n# 1) Swap the literal company name for the new placeholder variable.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Super Duper Inc.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"`$CURRENTUSERCOMPANY`$\"\n$find.Execute([ref]$find.Text, [ref]$true, $null, $null, $null, $null, $null, $null, $null, [ref]$find.Replacement.Text, 2)\n\n# 2) Remove the leftover \"_GoBack\" bookmark, if present.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
